$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48 (pushes existing rows 48:145 down to 49:146)
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly record
$ws.Cells.Item(48, 1).Value2 = 5
$ws.Cells.Item(48, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(48, 3).Value = "Maule"
$ws.Cells.Item(48, 4).Value2 = 45251
$ws.Cells.Item(48, 5).Value2 = 7
$ws.Cells.Item(48, 6).Value = "Fruta"
$ws.Cells.Item(48, 7).Value2 = 100101
$ws.Cells.Item(48, 8).Value = "Berries"
$ws.Cells.Item(48, 9).Value2 = 100101001
$ws.Cells.Item(48, 10).Value = "Arándano (blue)"
$ws.Cells.Item(48, 11).Value = "Sin especificar"
$ws.Cells.Item(48, 12).Value = "Primera"
$ws.Cells.Item(48, 13).Value2 = 100
$ws.Cells.Item(48, 14).Value2 = 6000
$ws.Cells.Item(48, 15).Value2 = 6000
$ws.Cells.Item(48, 16).Value2 = 6000
$ws.Cells.Item(48, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(48, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(48, 19).Value2 = 3000
$ws.Cells.Item(48, 20).Value2 = 2
